$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the first
#    (Heading1) paragraph: an empty run, a bold "Meta description" run, and
#    a regular run with the rest of the sentence.
# ---------------------------------------------------------------------------
$firstPara = $d.Paragraphs(1)
$firstPara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"
$metaPara.Range.Text = "Meta description: Explore solid payout potential in this Egyptian themed game with well-crafted graphics. Play 'Ancient Egypt' for free and check it out today."

$metaBoldRange = $metaPara.Range.Duplicate
$metaBoldRange.Find.Execute("Meta description", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$metaBoldRange.Bold = 1

# ---------------------------------------------------------------------------
# 2) Remove the trailing duplicate "Play 'Ancient Egypt' Free - ..." bold
#    title paragraph near the end of the document.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs($count - 1)
$dupTitlePara.Range.Delete()

# ---------------------------------------------------------------------------
# 3) Replace the text of the final (italic) paragraph with the new image
#    prompt copy, preserving its italic run formatting.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)

$oldItalicText = "Explore solid payout potential in this Egyptian themed game with well-crafted graphics. Play 'Ancient Egypt' for free and check it out today."
$newItalicText = "Please create a cartoon-style feature image for `"Ancient Egypt`" online slot game featuring a happy Maya warrior with glasses. The image should showcase the vibrant world of ancient Egypt, with rich colors and intricate details. The Maya warrior should be grinning, holding a golden scarab, and standing in front of the pyramids. The pyramids should feature a caricature style design, with bright colors and bold lines that complement the Maya warrior and add a playful touch to the image. The overall image should convey the excitement and thrill of playing an online slot game set in ancient Egypt and attract players who enjoy lively and entertaining games."

$italicTarget = $lastPara.Range.Duplicate
$italicTarget.Find.Execute($oldItalicText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$italicTarget.Text = $newItalicText
